# Update countries & provincias Spain
#
# 1) Reorder the country list: "Togo" moves up so it sits right after
#    "Bermudas" and before "Monaco" (previously it sat after "Cabo Verde").
#    Togo also gets refreshed case numbers. Monaco and Cabo Verde keep
#    their existing numbers but shift down one row.
# 2) Refresh the "Datos actualizados..." timestamp footer.
# 3) Refresh the case-count figures for several countries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Re-sequence Bermudas / Togo / Monaco / Cabo Verde / Zambia block ---
# Row 147 Bermudas stays put; Togo is inserted at row 148 (new data),
# Monaco drops to row 149 (same data it always had), Cabo Verde drops to
# row 150 (same data it always had); Zambia stays at row 151.

$ws.Cells.Item(148, 1).Value = "Togo"
$ws.Cells.Item(148, 2).Value = 96
$ws.Cells.Item(148, 3).Value = 6
$ws.Cells.Item(148, 4).Value = 62
$ws.Cells.Item(148, 5).Value = 28
$ws.Cells.Item(148, 6).Value = 0
$ws.Cells.Item(148, 7).Value = 0
$ws.Cells.Item(148, 8).Value = 6

$ws.Cells.Item(149, 1).Value = "Monaco"
$ws.Cells.Item(149, 2).Value = 94
$ws.Cells.Item(149, 3).Value = 0
$ws.Cells.Item(149, 4).Value = 41
$ws.Cells.Item(149, 5).Value = 49
$ws.Cells.Item(149, 6).Value = 2
$ws.Cells.Item(149, 7).Value = 0
$ws.Cells.Item(149, 8).Value = 4

$ws.Cells.Item(150, 1).Value = "Cabo Verde"
$ws.Cells.Item(150, 2).Value = 90
$ws.Cells.Item(150, 3).Value = 2
$ws.Cells.Item(150, 4).Value = 1
$ws.Cells.Item(150, 5).Value = 88
$ws.Cells.Item(150, 6).Value = 0
$ws.Cells.Item(150, 7).Value = 0
$ws.Cells.Item(150, 8).Value = 1

# --- 2) Update the footer timestamp text (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 25 de Abril de 2020 a las 15:52"

# --- 3) Update case-count figures for various countries ---

# Row 4 - Estados Unidos
$ws.Cells.Item(4, 2).Value = 926530
$ws.Cells.Item(4, 3).Value = 1298
$ws.Cells.Item(4, 5).Value = 763855
$ws.Cells.Item(4, 7).Value = 50
$ws.Cells.Item(4, 8).Value = 52243

# Row 5 - España
$ws.Cells.Item(5, 2).Value = 223759
$ws.Cells.Item(5, 3).Value = 3995
$ws.Cells.Item(5, 4).Value = 95708
$ws.Cells.Item(5, 5).Value = 105149
$ws.Cells.Item(5, 7).Value = 378
$ws.Cells.Item(5, 8).Value = 22902

# Row 43 - Noruega
$ws.Cells.Item(43, 2).Value = 7493
$ws.Cells.Item(43, 3).Value = 30
$ws.Cells.Item(43, 5).Value = 7260
$ws.Cells.Item(43, 7).Value = 2
$ws.Cells.Item(43, 8).Value = 201

# Row 81 - Republica de Macedonia
$ws.Cells.Item(81, 2).Value = 1367
$ws.Cells.Item(81, 3).Value = 41
$ws.Cells.Item(81, 4).Value = 374
$ws.Cells.Item(81, 5).Value = 934
$ws.Cells.Item(81, 7).Value = 2
$ws.Cells.Item(81, 8).Value = 59

# Row 107 - Georgia
$ws.Cells.Item(107, 4).Value = 139
$ws.Cells.Item(107, 5).Value = 312

# Row 109 - Jordania
$ws.Cells.Item(109, 2).Value = 444
$ws.Cells.Item(109, 3).Value = 3
$ws.Cells.Item(109, 4).Value = 332
$ws.Cells.Item(109, 5).Value = 105

# Row 110 - Sri Lanka
$ws.Cells.Item(110, 2).Value = 435
$ws.Cells.Item(110, 3).Value = 18
$ws.Cells.Item(110, 4).Value = 118

# Row 116 - Kenia
$ws.Cells.Item(116, 4).Value = 98
$ws.Cells.Item(116, 5).Value = 231
